$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I0, IF) styled like the existing header cells (bold, bordered, centered)
$ws.Range("I1").Value = "I0"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.LineStyle = 1

$ws.Range("J1").Value = "IF"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

# New data values for columns I (I0) and J (IF)
$iValues = @(9, 5, 7, 6, 3, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 3)
$jValues = @(9, 6, 9, 7, 6, 5, 4, 5, 5, 5, 5, 6, 5, 4, 3, 2, 3)

for ($r = 2; $r -le 18; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
